# Updated cryptos list on Sat Apr 27 11:51:43 UTC 2024 with GitHub Actions
#
# Sets a cell's value while keeping it a genuine text string (matching the
# source data, which stores every Price/Volume cell as text even when it
# looks numeric, e.g. "589.67"). Plain "$ws.Range($addr).Value = $val" would
# let Excel auto-coerce a numeric-looking string into a real number, so we
# briefly mark the cell as Text ("@") before writing, then restore the
# "Normal" style so no stray number-format is left behind on the cell.
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.955.16"
$ws.Range("E2").Value = "  -2.04%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.130.08"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "589.15"
$ws.Range("E5").Value = "  -2.21%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "136.13"
$ws.Range("E6").Value = "  -4.86%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.121.89"
$ws.Range("E8").Value = "  -0.17%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -1.69%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.90%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -3.11%  "

# Row 12 - Cardano
Set-TextValue $ws "D12" "0.453"
$ws.Range("E12").Value = "  -3.42%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -5.39%  "

# Row 14 - Avalanche
Set-TextValue $ws "D14" "33.91"
$ws.Range("E14").Value = "  -3.71%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.638.76"
$ws.Range("E15").Value = "  -0.15%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +1.38%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "62.994.24"
$ws.Range("E17").Value = "  -1.78%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.123.72"
$ws.Range("E18").Value = "  +0.11%  "

# Row 19 - Polkadot
Set-TextValue $ws "D19" "6.61"
$ws.Range("E19").Value = "  -4.10%  "

# Row 20 - BitcoinCash
Set-TextValue $ws "D20" "469.90"
$ws.Range("E20").Value = "  -2.52%  "

# Row 21 - Chainlink
Set-TextValue $ws "D21" "14.05"
$ws.Range("E21").Value = "  -3.97%  "

# Row 22 - Polygon
Set-TextValue $ws "D22" "0.695"
$ws.Range("E22").Value = "  -2.20%  "

# Row 23 - Uniswap
Set-TextValue $ws "D23" "7.63"
$ws.Range("E23").Value = "  -0.81%  "

# Row 24 - Litecoin
Set-TextValue $ws "D24" "85.24"
$ws.Range("E24").Value = "  -0.38%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue $ws "D25" "12.89"
$ws.Range("E25").Value = "  -4.15%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  -1.86%  "

# Row 28 - RenderToken
Set-TextValue $ws "D28" "7.86"
$ws.Range("E28").Value = "  -5.86%  "

# Row 29 - ImmutableX
Set-TextValue $ws "D29" "2.08"
$ws.Range("E29").Value = "  +1.52%  "

# Row 30 - NEARProtocol
Set-TextValue $ws "D30" "6.82"
$ws.Range("E30").Value = "  -5.00%  "

# Row 31 - FirstDigitalUSD
$ws.Range("E31").Value = "  +0.02%  "

# Row 32 - EthereumClassic
Set-TextValue $ws "D32" "26.52"
$ws.Range("E32").Value = "  -1.06%  "

# Row 33 - Hedera
Set-TextValue $ws "D33" "0.107"
$ws.Range("E33").Value = "  -4.61%  "

# Row 34 - Stacks
$ws.Range("E34").Value = "  -4.69%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -3.33%  "

# Rows 36/37 - Filecoin and OKB swap rank positions (OKB now ranks above Filecoin)
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D36" "51.94"
$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D37" "5.73"
$ws.Range("E37").Value = "  -4.18%  "

# Row 38 - PEPE
$ws.Range("E38").Value = "  -12.19%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -2.41%  "

# Row 40 - Bittensor
Set-TextValue $ws "D40" "415.01"
$ws.Range("E40").Value = "  -7.06%  "

# Row 41 - Maker
$ws.Range("D41").Value = "2.913.79"
$ws.Range("E41").Value = "  +1.96%  "

# Row 42 - Cosmos
Set-TextValue $ws "D42" "8.16"
$ws.Range("E42").Value = "  -0.74%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -11.34%  "

# Row 44 - Kaspa
$ws.Range("E44").Value = "  -6.35%  "

# Row 45 - TheGraph
Set-TextValue $ws "D45" "0.258"
$ws.Range("E45").Value = "  -0.77%  "

# Row 47 - Fetch.AI
$ws.Range("E47").Value = "  -5.84%  "

# Row 48 - InjectiveProtocol
Set-TextValue $ws "D48" "25.23"
$ws.Range("E48").Value = "  -3.06%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -0.98%  "

# Row 50 - ThetaToken
Set-TextValue $ws "D50" "2.24"
$ws.Range("E50").Value = "  -8.06%  "

# Row 51 - Monero
Set-TextValue $ws "D51" "120.23"
$ws.Range("E51").Value = "  -0.51%  "
